$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.69225069466188
$ws.Range("D2").Value = 10.0966868703377
$ws.Range("E2").Value = 0.0969933546772152

$ws.Range("C3").Value = 2.69500226854229
$ws.Range("D3").Value = 19.1264517407852
$ws.Range("E3").Value = 0.302986326862212

$ws.Range("C4").Value = 1.44306057544509
$ws.Range("D4").Value = 10.0970941454359
$ws.Range("E4").Value = 0.239852520712065

$ws.Range("C5").Value = 1.94050973213217
$ws.Range("D5").Value = 13.3279090222345
$ws.Range("E5").Value = 0.106344302128826

$ws.Range("C6").Value = 0.564853880819194
$ws.Range("D6").Value = 3.84764516101637
$ws.Range("E6").Value = 0.0370508657663885
$ws.Range("F6").Value = 0.0016
$ws.Range("G6").Value = 0.0096

$ws.Range("C7").Value = 2.11109306201665
$ws.Range("D7").Value = 20.5941144479453
$ws.Range("E7").Value = 0.291725657429005
